# "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 5 (Damasco, variety "Dina",
# Provincia de Quillota), pushing all the existing data rows (previously 5..45)
# down to 6..46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 5, shifting rows 5-45 down to 6-46.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 'Vega Modelo de Temuco'
$ws.Range("C5").Value = 'La Araucanía'
$ws.Range("D5").Value = 44550
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = 'Frutos de hueso (carozo)'
$ws.Range("I5").Value = 100103003
$ws.Range("J5").Value = 'Damasco'
$ws.Range("K5").Value = 'Dina'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 21000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Provincia de Quillota'
$ws.Range("S5").Value = 1167
$ws.Range("T5").Value = 18
